# Settings save button fixed
#
# Adds the new "AB" location (for company "NHS") and the two asset-type
# rows that go with it:
#   - Weir / BC / NHS / #00a80e
#   - Cableway / AB / NHS / #70ffdb

$wb = $excel.ActiveWorkbook

# --- Locations sheet: add new location "AB" for company "NHS" ---
$locations = $wb.Worksheets.Item("Locations")
$locRow = $locations.UsedRange.Rows.Count + 1
$locations.Cells.Item($locRow, 1).Value = "AB"
$locations.Cells.Item($locRow, 2).Value = "NHS"

# --- AssetTypes sheet: add the two new asset type rows ---
$assetTypes = $wb.Worksheets.Item("AssetTypes")
$atRow = $assetTypes.UsedRange.Rows.Count + 1

$assetTypes.Cells.Item($atRow, 1).Value = "Weir"
$assetTypes.Cells.Item($atRow, 2).Value = "BC"
$assetTypes.Cells.Item($atRow, 3).Value = "NHS"
$assetTypes.Cells.Item($atRow, 4).Value = "#00a80e"
$atRow = $atRow + 1

$assetTypes.Cells.Item($atRow, 1).Value = "Cableway"
$assetTypes.Cells.Item($atRow, 2).Value = "AB"
$assetTypes.Cells.Item($atRow, 3).Value = "NHS"
$assetTypes.Cells.Item($atRow, 4).Value = "#70ffdb"
